$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. TC_002 row (B3): EXECUTE flag flips from YES to NO
$ws.Range("B3").Value = "NO"

# 2. Fix up C16's leftover "header-ish" bold style so it matches the plain
#    data-row style used by its neighbours before we drop a value into it.
$ws.Range("D16").Copy()
$ws.Range("C16").PasteSpecial(-4122)

# 3. New API test case, split across 4 rows (one per HTTP verb / call).
#    Values are entered column-by-column (as the original author did) so the
#    shared-string table grows in the same order.
$ws.Range("C16:C19").Value = "API"
$ws.Range("D16:D19").Value = "ApiCalls"

$ws.Range("E16").Style = "Normal"
$ws.Range("E16").Value = "getSingleUserCall"
$ws.Range("E16").VerticalAlignment = -4108

$ws.Range("E17").Value = "postCreateUserCall"
$ws.Range("E17").VerticalAlignment = -4108

$ws.Range("A16:A19").Value = "TC_014_GET_CALL"
$ws.Range("B16:B19").Value = "YES"

$ws.Range("E18").Value = "putUpdateUserCall"
$ws.Range("E18").VerticalAlignment = -4108

$ws.Range("E19").Value = "deleteUserCall"
$ws.Range("E19").VerticalAlignment = -4108

# 4. Row 16 loses its old oversized "section title" height now that it holds
#    a normal data row
$ws.Rows.Item(16).AutoFit()

# 5. Scroll / selection state left by the editor
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("D22").Select()
